$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Resueltos (B) / Pendientes (C) values for several departments ---
# Row 2 - Tesoreria
$ws.Range("B2").Value = 93
# Row 3 - Administracion
$ws.Range("B3").Value = 20
# Row 4 - Consejeria_Juridica
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 18
# Row 6 - Desarrollo_Economico_Turistico_y_Artesanal
$ws.Range("B6").Value = 10
# Row 7 - Desarrollo_Social
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 5
# Row 8 - Seguridad_Publica
$ws.Range("B8").Value = 40
# Row 9 - Desarrollo_Urbano_y_Metropolitano
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 2
# Row 12 - Medio_Ambiente
$ws.Range("B12").Value = 22
# Row 13 - Gobierno_por_Resultados
$ws.Range("C13").Value = 2
# Row 14 - Igualdad_de_Genero
$ws.Range("B14").Value = 9
# Row 15 - Educacion
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 0
# Row 16 - Contraloria_Municipal
$ws.Range("B16").Value = 26
$ws.Range("C16").Value = 7
# Row 17 - OPDAPAS
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 5
# Row 18 - IMCUFIDEM
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = 0
# Row 19 - SMDIF
$ws.Range("C19").Value = 1
# Row 24 - Proteccion_Civil_y_Bomberos
$ws.Range("B24").Value = 19
# Row 25 - Gobierno_digital
$ws.Range("C25").Value = 2
# Row 27 - Obras_Publicas
$ws.Range("C27").Value = 6
# Row 28 - Defensora_Municipal_de_los_derechos_humanos
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 2

# --- E column: replace the static "total tareas" (366) with a live formula
#     pointing at the grand total in D29, and add the same to E1 ---
# (set cell-by-cell rather than over the whole range at once so Excel does
#  not turn this into a shared formula whose relative reference would shift
#  row-by-row)
$ws.Range("E1").Formula = "=D29"
for ($r = 2; $r -le 28; $r++) {
    $ws.Range("E$r").Formula = "=D29"
}

# --- Move the active selection from I17 to H19 ---
$ws.Range("H19").Select()
